$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = 111909536
$ws.Range("B14").Value = 77388
$ws.Range("E14").Value = 6446
$ws.Range("F14").Value = "Kolflarnlav"
$ws.Range("G14").Value = "Carbonicola anthracophila"
$ws.Range("H14").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("P14").Value = "Fläcksberget, Hjd"
$ws.Range("Q14").Value = 467891
$ws.Range("R14").Value = 6875425

# Row 15
$ws.Range("A15").Value = 111908364
$ws.Range("B15").Value = 90794
$ws.Range("E15").Value = 4362
$ws.Range("F15").Value = "Blå taggsvamp"
$ws.Range("G15").Value = "Hydnellum caeruleum"
$ws.Range("H15").Value = "(Hornem.) P.Karst."
$ws.Range("P15").Value = "Gröbäcken, Hjd"
$ws.Range("Q15").Value = 467724
$ws.Range("R15").Value = 6874811

# Row 16
$ws.Range("A16").Value = 111909174
$ws.Range("B16").Value = 77388
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 6446
$ws.Range("F16").Value = "Kolflarnlav"
$ws.Range("G16").Value = "Carbonicola anthracophila"
$ws.Range("H16").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("I16").ClearContents()
$ws.Range("J16").ClearContents()
$ws.Range("Q16").Value = 467989
$ws.Range("R16").Value = 6875353
$ws.Range("AC16").ClearContents()

# Row 17
$ws.Range("B17").Value = 96720

# Row 18
$ws.Range("A18").Value = 111908700
$ws.Range("B18").Value = 90099
$ws.Range("D18").Value = "VU"
$ws.Range("E18").Value = 760
$ws.Range("F18").Value = "Doftticka"
$ws.Range("G18").Value = "Haploporus odorus"
$ws.Range("H18").Value = "(Sommerf.) Bondartsev & Singer"
$ws.Range("I18").Value = "6"
$ws.Range("J18").Value = "fruktkroppar"
$ws.Range("Q18").Value = 467922
$ws.Range("R18").Value = 6875307
$ws.Range("AC18").Value = "Förekomst av doftticka i avverkningsanmält område."

# Row 19
$ws.Range("A19").Value = 111909766
$ws.Range("B19").Value = 89317
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 3215
$ws.Range("F19").Value = "Rödgul trumpetsvamp"
$ws.Range("G19").Value = "Craterellus lutescens"
$ws.Range("H19").Value = "(Fr.) Fr."
$ws.Range("Q19").Value = 467757
$ws.Range("R19").Value = 6875470

# Row 20
$ws.Range("B20").Value = 77388

# Row 21
$ws.Range("B21").Value = 90812

# Row 22
$ws.Range("A22").Value = 112014300
$ws.Range("B22").Value = 90823
$ws.Range("E22").Value = 5966
$ws.Range("F22").Value = "Motaggsvamp"
$ws.Range("G22").Value = "Sarcodon squamosus"
$ws.Range("H22").Value = "(Schaeff.) Quél."
$ws.Range("Q22").Value = 467415
$ws.Range("R22").Value = 6875287

# Row 23
$ws.Range("A23").Value = 112014923
$ws.Range("B23").Value = 90823
$ws.Range("Q23").Value = 467413
$ws.Range("R23").Value = 6875234

# Row 24
$ws.Range("A24").Value = 112014177
$ws.Range("B24").Value = 90823
$ws.Range("E24").Value = 5966
$ws.Range("F24").Value = "Motaggsvamp"
$ws.Range("G24").Value = "Sarcodon squamosus"
$ws.Range("H24").Value = "(Schaeff.) Quél."

# Row 25
$ws.Range("A25").Value = 112015011
$ws.Range("B25").Value = 90792
$ws.Range("Q25").Value = 467390
$ws.Range("R25").Value = 6875328

# Row 26
$ws.Range("A26").Value = 112014423
$ws.Range("B26").Value = 90792
$ws.Range("E26").Value = 4361
$ws.Range("F26").Value = "Orange taggsvamp"
$ws.Range("G26").Value = "Hydnellum aurantiacum"
$ws.Range("H26").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q26").Value = 467430
$ws.Range("R26").Value = 6875238

# Row 27
$ws.Range("B27").Value = 90816

# Row 28
$ws.Range("A28").Value = 112014208
$ws.Range("B28").Value = 90792
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 4361
$ws.Range("F28").Value = "Orange taggsvamp"
$ws.Range("G28").Value = "Hydnellum aurantiacum"
$ws.Range("H28").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q28").Value = 467418
$ws.Range("R28").Value = 6875313

# Row 29
$ws.Range("A29").Value = 112014142
$ws.Range("B29").Value = 90800
$ws.Range("D29").Value = "LC"
$ws.Range("E29").Value = 4364
$ws.Range("F29").Value = "Dropptaggsvamp"
$ws.Range("G29").Value = "Hydnellum ferrugineum"
$ws.Range("H29").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q29").Value = 467443
$ws.Range("R29").Value = 6875337
